$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-10) have been reordered by ascending date,
# with columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion), R (Origen),
# S (Precio $/Kg) and T (Kg / unidad) following the row they belong to.
# Columns A, B, C, E, F, G, H, I, J, K, L stay identical for every row.

$rows = @{
    2  = @{ D = 44540; M = 240; N = 3500; O = 3800; P = 3650; Q = "`$/bandeja 2 kilos"; R = "Región del Maule";       S = 1825; T = 2 }
    3  = @{ D = 44181; M = 65;  N = 3600; O = 3800; P = 3692; Q = "`$/bandeja 2 kilos"; R = "Provincia de Diguillín"; S = 1846; T = 2 }
    4  = @{ D = 44181; M = 80;  N = 1800; O = 2000; P = 1875; Q = "`$/envase 1 kilo";   R = "Provincia de Diguillín"; S = 1875; T = 1 }
    5  = @{ D = 44174; M = 150; N = 3700; O = 3800; P = 3747; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";   S = 1874; T = 2 }
    6  = @{ D = 44187; M = 80;  N = 2800; O = 3000; P = 2900; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";   S = 1450; T = 2 }
    7  = @{ D = 44187; M = 65;  N = 1400; O = 1500; P = 1446; Q = "`$/envase 1 kilo";   R = "Provincia de Diguillín"; S = 1446; T = 1 }
    8  = @{ D = 44596; M = 120; N = 2500; O = 2700; P = 2600; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";   S = 1300; T = 2 }
    9  = @{ D = 44539; M = 200; N = 3800; O = 4000; P = 3900; Q = "`$/bandeja 2 kilos"; R = "Región del Maule";       S = 1950; T = 2 }
    10 = @{ D = 44594; M = 120; N = 2500; O = 2800; P = 2650; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";   S = 1325; T = 2 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $data.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $data.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $data.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $data.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $data.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $data.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $data.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $data.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $data.T   # T: Kg / unidad
}
